$wb = $excel.ActiveWorkbook

# Existing "Productos" sheet - the new cart sheet mirrors its product list.
$productos = $wb.Worksheets.Item("Productos")

# Add the new "ProductosCarrito" worksheet right after "Productos".
$carrito = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $productos)
$carrito.Name = "ProductosCarrito"

# Populate it with the same product names as "Productos".
$carrito.Range("A1").Value = "Nombre"
$carrito.Range("A2").Value = "Iphone"
$carrito.Range("A3").Value = "MacBook"
$carrito.Range("A4").Value = "Samsung Galaxy Tab"

# Match the cell style used on the "Productos" sheet.
$carrito.Range("A1:A4").Style = $productos.Range("A1:A4").Style

# Set the selection on the new sheet and make it the active tab.
$carrito.Range("J11").Select()
$carrito.Activate()

# Update the selection on the "Productos" sheet (no longer the active tab).
$productos.Range("A1:A4").Select()
